$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the Cypher query text in C2:C4 (shared by all three rows) ---
$newQuery = "MATCH (p:program)<--(s:study)<-[*]-(c:case)<--(demo:demographic)
OPTIONAL MATCH (samp:sample)-->(c)
OPTIONAL MATCH (diag:diagnosis)-->(c)
OPTIONAL MATCH (f:file)-[*]->(c)
OPTIONAL MATCH (sf:file)-->(s)
WITH DISTINCT f, sf, samp AS samp, c, demo, diag, s, p
WHERE demo.breed IN ['German Shorthaired Pointer']
RETURN  
    count(distinct p) AS Programs,
    count(distinct s) AS Studies,
    count(distinct c) AS Cases,
    count(distinct samp) AS Samples,
    count(distinct f) AS ``Case Files``,
    count(distinct sf) AS ``Study Files``"

$ws.Range("C2").Value = $newQuery
$ws.Range("C3").Value = $newQuery
$ws.Range("C4").Value = $newQuery

# --- Adjust row heights for rows 2-4 (previously pinned at the 409.6 max) ---
$ws.Rows(2).RowHeight = 244.8
$ws.Rows(3).RowHeight = 230.4
$ws.Rows(4).RowHeight = 216

# --- Update zoom level on the active window ---
$excel.ActiveWindow.Zoom = 115

# --- Update the selected range ---
$ws.Range("B4:B5").Select() | Out-Null
